$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2570887976341871
$ws.Range("C2").Value = 5951653278551590

$ws.Range("B3").Value = 2535312476896360
$ws.Range("C3").Value = 5869295631216111

$ws.Range("B4").Value = 2526211061933456
$ws.Range("C4").Value = 5848226765878424

$ws.Range("B5").Value = 1835035606876774
$ws.Range("C5").Value = 4248145187274770

$ws.Range("B6").Value = 1125092206874587
$ws.Range("C6").Value = 2604612672976170

$ws.Range("B7").Value = 1051984338854707
$ws.Range("C7").Value = 2435367019846647

$ws.Range("B8").Value = 790743204785341.4
$ws.Range("C8").Value = 1830587821417162

$ws.Range("B9").Value = 128928077553395.9
$ws.Range("C9").Value = 298471340025910.3

$ws.Range("B10").Value = 141059503.6121289
$ws.Range("C10").Value = 326555671.994087
